$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.55
$ws.Range("G2").Value = 1.56
$ws.Range("H2").Value = 6.2
$ws.Range("I2").Value = 6.6
$ws.Range("J2").Value = 4.9
$ws.Range("K2").Value = 5.1
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.19
$ws.Range("P2").Value = 2.64
$ws.Range("Q2").Value = 1.59
$ws.Range("R2").Value = 1.64
$ws.Range("S2").Value = 2.48
$ws.Range("T2").Value = 1.73
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.18
$ws.Range("W2").Value = 2.78
$ws.Range("X2").Value = 27
$ws.Range("Y2").Value = 30
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 160
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 85
$ws.Range("AF2").Value = 10.5
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 65
$ws.Range("AJ2").Value = 14.5
$ws.Range("AK2").Value = 14
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 6.2
$ws.Range("AO2").Value = 70
$ws.Range("F3").Value = 25
$ws.Range("G3").Value = 34
$ws.Range("H3").Value = 1.09
$ws.Range("I3").Value = 1.1
$ws.Range("K3").Value = 19
$ws.Range("L3").Value = 1.11
$ws.Range("N3").Value = 16
$ws.Range("O3").Value = 1.05
$ws.Range("P3").Value = 6.8
$ws.Range("R3").Value = 3.4
$ws.Range("T3").Value = 1.81
$ws.Range("U3").Value = 1.97
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 1.03
$ws.Range("Y3").Value = 34
$ws.Range("AC3").Value = 55
$ws.Range("AD3").Value = 21
$ws.Range("AH3").Value = 990
$ws.Range("AI3").Value = 38
$ws.Range("AO3").Value = 1.88
$ws.Range("F4").Value = 2.3
$ws.Range("G4").Value = 2.48
$ws.Range("H4").Value = 3.85
$ws.Range("I4").Value = 4.5
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1.73
$ws.Range("M4").Value = 1.18
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.75
$ws.Range("P4").Value = 1.4
$ws.Range("Q4").Value = 3.4
$ws.Range("R4").Value = 1.12
$ws.Range("X4").Value = 6.6
$ws.Range("Y4").Value = 9.8
$ws.Range("Z4").Value = 80
$ws.Range("AB4").Value = 6.4
$ws.Range("AC4").Value = 7.4
$ws.Range("AD4").Value = 36
$ws.Range("AF4").Value = 13.5
$ws.Range("AG4").Value = 14.5
$ws.Range("AH4").Value = 65
$ws.Range("AJ4").Value = 48
$ws.Range("AN4").Value = 120
$ws.Range("F5").Value = 2.36
$ws.Range("G5").Value = 2.46
$ws.Range("H5").Value = 2.94
$ws.Range("I5").Value = 3
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.29
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 2.68
$ws.Range("Q5").Value = 1.58
$ws.Range("R5").Value = 1.68
$ws.Range("S5").Value = 2.4
$ws.Range("T5").Value = 1.51
$ws.Range("U5").Value = 2.68
$ws.Range("V5").Value = 1.5
$ws.Range("W5").Value = 1.69
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 18.5
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 50
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 9.6
$ws.Range("AD5").Value = 14
$ws.Range("AE5").Value = 28
$ws.Range("AF5").Value = 19.5
$ws.Range("AG5").Value = 12.5
$ws.Range("AI5").Value = 32
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 23
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 55
$ws.Range("AN5").Value = 12
$ws.Range("AO5").Value = 17
$ws.Range("F6").Value = 3.2
$ws.Range("G6").Value = 3.85
$ws.Range("H6").Value = 2.64
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 2.64
$ws.Range("P6").Value = 1.43
$ws.Range("U6").Value = 1.66
$ws.Range("V6").Value = 1.51
$ws.Range("W6").Value = 1.37
$ws.Range("X6").Value = 11
$ws.Range("AA6").Value = 220
$ws.Range("AJ6").Value = 1000
$ws.Range("G7").Value = 1.47
$ws.Range("H7").Value = 8.6
$ws.Range("I7").Value = 9.4
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.5
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 2.2
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.45
$ws.Range("S7").Value = 2.9
$ws.Range("T7").Value = 1.96
$ws.Range("U7").Value = 1.89
$ws.Range("V7").Value = 1.12
$ws.Range("W7").Value = 3.1
$ws.Range("X7").Value = 18.5
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 190
$ws.Range("AA7").Value = 300
$ws.Range("AC7").Value = 11.5
$ws.Range("AG7").Value = 10
$ws.Range("AH7").Value = 26
$ws.Range("AJ7").Value = 12.5
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 85
$ws.Range("AM7").Value = 150
$ws.Range("AN7").Value = 7.2
$ws.Range("AO7").Value = 1000
$ws.Range("F8").Value = 2.74
$ws.Range("G8").Value = 2.9
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 2.88
$ws.Range("N8").Value = 3.55
$ws.Range("O8").Value = 1.34
$ws.Range("Q8").Value = 2.04
$ws.Range("T8").Value = 1.77
$ws.Range("V8").Value = 1.53
$ws.Range("Z8").Value = 21
$ws.Range("AE8").Value = 95
$ws.Range("AH8").Value = 21
$ws.Range("AO8").Value = 600
$ws.Range("F9").Value = 1.58
$ws.Range("G9").Value = 1.67
$ws.Range("H9").Value = 9.2
$ws.Range("K9").Value = 3.8
$ws.Range("L9").Value = 1.71
$ws.Range("M9").Value = 1.16
$ws.Range("N9").Value = 2.22
$ws.Range("O9").Value = 1.7
$ws.Range("P9").Value = 1.4
$ws.Range("R9").Value = 1.12
$ws.Range("S9").Value = 7.4
$ws.Range("T9").Value = 3
$ws.Range("U9").Value = 1.44
$ws.Range("W9").Value = 2.48
$ws.Range("Y9").Value = 60
$ws.Range("AB9").Value = 4.8
$ws.Range("AC9").Value = 14
$ws.Range("AF9").Value = 13.5
$ws.Range("F10").Value = 1.67
$ws.Range("G10").Value = 1.68
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 7.2
$ws.Range("J10").Value = 3.8
$ws.Range("K10").Value = 3.85
$ws.Range("L10").Value = 1.54
$ws.Range("O10").Value = 1.46
$ws.Range("P10").Value = 1.69
$ws.Range("Q10").Value = 2.4
$ws.Range("R10").Value = 1.25
$ws.Range("S10").Value = 4.8
$ws.Range("U10").Value = 1.71
$ws.Range("W10").Value = 2.46
$ws.Range("X10").Value = 10.5
$ws.Range("AA10").Value = 250
$ws.Range("AD10").Value = 27
$ws.Range("AE10").Value = 130
$ws.Range("AK10").Value = 21
$ws.Range("AM10").Value = 220
$ws.Range("AN10").Value = 14
$ws.Range("H11").Value = 46
$ws.Range("Q11").Value = 1.44
$ws.Range("R11").Value = 1.85
$ws.Range("S11").Value = 2.12
$ws.Range("T11").Value = 2.74
$ws.Range("U11").Value = 1.53
$ws.Range("X11").Value = 60
$ws.Range("Z11").Value = 560
$ws.Range("AB11").Value = 12
$ws.Range("AC11").Value = 27
$ws.Range("AG11").Value = 15.5
$ws.Range("AH11").Value = 990
$ws.Range("AK11").Value = 15.5
$ws.Range("AM11").Value = 500
$ws.Range("AN11").Value = 2.88
$ws.Range("M12").Value = 1.02
$ws.Range("O12").Value = 1.02
$ws.Range("Q12").Value = 1.06
